$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1378.3334
$ws.Range("I31").Value = 1072.1428
$ws.Range("K31").Value = 3216.4284
$ws.Range("M31").Value = -2986.4284
$ws.Range("H33").Value = 145.33333
$ws.Range("I33").Value = 150.71428
$ws.Range("J33").Value = 107.666664
$ws.Range("K33").Value = 150.71428
$ws.Range("L33").Value = 107.666664
$ws.Range("M33").Value = 78.28572
$ws.Range("N33").Value = -565.666664
$ws.Range("H40").Value = 2148.3
$ws.Range("I40").Value = 1928.4286
$ws.Range("J40").Value = 2215.2173
$ws.Range("K40").Value = 1928.4286
$ws.Range("L40").Value = 2215.2173
$ws.Range("M40").Value = -1753.4286
$ws.Range("N40").Value = -2565.2173
$ws.Range("H62").Value = 97109.45
$ws.Range("I62").Value = 128975.5
$ws.Range("K62").Value = 128975.5
$ws.Range("M62").Value = -128351.5
$ws.Range("H64").Value = 3535.5
$ws.Range("I64").Value = 3350.2
$ws.Range("J64").Value = 3651.3125
$ws.Range("K64").Value = 3350.2
$ws.Range("L64").Value = 3651.3125
$ws.Range("M64").Value = -3102.2
$ws.Range("N64").Value = -4147.3125
$ws.Range("H65").Value = 97109.45
$ws.Range("I65").Value = 128975.5
$ws.Range("K65").Value = 644877.5
$ws.Range("M65").Value = -641757.5
$ws.Range("H67").Value = 3535.5
$ws.Range("I67").Value = 3350.2
$ws.Range("J67").Value = 3651.3125
$ws.Range("K67").Value = 3350.2
$ws.Range("L67").Value = 3651.3125
$ws.Range("M67").Value = -2492.2
$ws.Range("N67").Value = -5367.3125
$ws.Range("H125").Value = 1220.8
$ws.Range("I125").Value = 1242.8572
$ws.Range("J125").Value = 1201.5
$ws.Range("K125").Value = 11185.7148
$ws.Range("L125").Value = 10813.5
$ws.Range("M125").Value = -8725.7148
$ws.Range("N125").Value = -15733.5
$ws.Range("H132").Value = 1750.125
$ws.Range("I132").Value = 1708.7037
$ws.Range("J132").Value = 2229.4285
$ws.Range("K132").Value = 5126.1111
$ws.Range("L132").Value = 6688.2855
$ws.Range("M132").Value = -2596.1111
$ws.Range("N132").Value = -11748.2855
$ws.Range("H135").Value = 989.5682
$ws.Range("I135").Value = 621.0294
$ws.Range("K135").Value = 5589.2646
$ws.Range("M135").Value = -3054.2646
$ws.Range("H137").Value = 755.4074000000001
$ws.Range("I137").Value = 738.3077
$ws.Range("K137").Value = 2214.9231
$ws.Range("M137").Value = 335.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3804984.2
$ws.Range("I32").Value = 4648624
$ws.Range("J32").Value = 8603.9375
$ws.Range("K32").Value = 4648624
$ws.Range("L32").Value = 8603.9375
$ws.Range("M32").Value = -4648337
$ws.Range("N32").Value = -9177.9375
$ws.Range("H45").Value = 1533.3572
$ws.Range("I45").Value = 1252.4546
$ws.Range("J45").Value = 2563.3333
$ws.Range("K45").Value = 1252.4546
$ws.Range("L45").Value = 2563.3333
$ws.Range("M45").Value = -875.4546
$ws.Range("N45").Value = -3317.3333
$ws.Range("H110").Value = 481.65518
$ws.Range("I110").Value = 428.08694
$ws.Range("J110").Value = 687
$ws.Range("K110").Value = 428.08694
$ws.Range("L110").Value = 687
$ws.Range("M110").Value = 1616.91306
$ws.Range("N110").Value = -4777
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800
$ws.Range("H132").Value = 1895.9445
$ws.Range("I132").Value = 1101.75
$ws.Range("K132").Value = 3305.25
$ws.Range("M132").Value = -775.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 29453.81
$ws.Range("I134").Value = 2147.92
$ws.Range("J134").Value = 86341.086
$ws.Range("K134").Value = 6443.76
$ws.Range("L134").Value = 259023.258
$ws.Range("M134").Value = -3908.76
$ws.Range("N134").Value = -264093.258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1470.2195
$ws.Range("I31").Value = 1456.625
$ws.Range("J31").Value = 2014
$ws.Range("K31").Value = 1456.625
$ws.Range("L31").Value = 2014
$ws.Range("M31").Value = -1161.625
$ws.Range("N31").Value = -2604
$ws.Range("H34").Value = 1470.2195
$ws.Range("I34").Value = 1456.625
$ws.Range("J34").Value = 2014
$ws.Range("K34").Value = 1456.625
$ws.Range("L34").Value = 2014
$ws.Range("M34").Value = -1254.625
$ws.Range("N34").Value = -2418
$ws.Range("H99").Value = 2548.1936
$ws.Range("I99").Value = 2215.7896
$ws.Range("J99").Value = 3074.5
$ws.Range("K99").Value = 2215.7896
$ws.Range("L99").Value = 3074.5
$ws.Range("M99").Value = -717.7896000000001
$ws.Range("N99").Value = -6070.5
$ws.Range("H107").Value = 264.0645
$ws.Range("I107").Value = 187.625
$ws.Range("J107").Value = 345.6
$ws.Range("K107").Value = 187.625
$ws.Range("L107").Value = 345.6
$ws.Range("M107").Value = 1732.375
$ws.Range("N107").Value = -4185.6
$ws.Range("H126").Value = 2548.1936
$ws.Range("I126").Value = 2215.7896
$ws.Range("J126").Value = 3074.5
$ws.Range("K126").Value = 6647.3688
$ws.Range("L126").Value = 9223.5
$ws.Range("M126").Value = -4177.3688
$ws.Range("N126").Value = -14163.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4826371.5
$ws.Range("I22").Value = 7239323.5
$ws.Range("J22").Value = 466.66666
$ws.Range("K22").Value = 21717970.5
$ws.Range("L22").Value = 1399.99998
$ws.Range("M22").Value = -21717801.5
$ws.Range("N22").Value = -1737.99998
$ws.Range("H27").Value = 4826371.5
$ws.Range("I27").Value = 7239323.5
$ws.Range("J27").Value = 466.66666
$ws.Range("K27").Value = 21717970.5
$ws.Range("L27").Value = 1399.99998
$ws.Range("M27").Value = -21717868.5
$ws.Range("N27").Value = -1603.99998
$ws.Range("H42").Value = 1875
$ws.Range("J42").Value = 1875
$ws.Range("L42").Value = 5625
$ws.Range("N42").Value = -6693
$ws.Range("H56").Value = 2943.7778
$ws.Range("I56").Value = 2943.7778
$ws.Range("K56").Value = 2943.7778
$ws.Range("M56").Value = -2413.7778
$ws.Range("H92").Value = 1560.4
$ws.Range("I92").Value = 267.33334
$ws.Range("J92").Value = 3500
$ws.Range("K92").Value = 802.0000200000001
$ws.Range("L92").Value = 10500
$ws.Range("M92").Value = 445.9999799999999
$ws.Range("N92").Value = -12996
$ws.Range("H114").Value = 1315.1428
$ws.Range("I114").Value = 681.2857
$ws.Range("J114").Value = 1949
$ws.Range("K114").Value = 2043.8571
$ws.Range("L114").Value = 5847
$ws.Range("M114").Value = 1210.1429
$ws.Range("N114").Value = -12355
$ws.Range("H116").Value = 127316
$ws.Range("I116").Value = 1305.8
$ws.Range("J116").Value = 337333
$ws.Range("K116").Value = 3917.4
$ws.Range("L116").Value = 1011999
$ws.Range("M116").Value = -475.3999999999996
$ws.Range("N116").Value = -1018883
$ws.Range("H131").Value = 779.11
$ws.Range("I131").Value = 355.53333
$ws.Range("J131").Value = 853.8588
$ws.Range("K131").Value = 1066.59999
$ws.Range("L131").Value = 2561.5764
$ws.Range("M131").Value = 3973.40001
$ws.Range("N131").Value = -12641.5764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.5
$ws.Range("I2").Value = 110.75
$ws.Range("J2").Value = 101
$ws.Range("K2").Value = 110.75
$ws.Range("L2").Value = 101
$ws.Range("M2").Value = 2.25
$ws.Range("N2").Value = -327
$ws.Range("H70").Value = 4284.9443
$ws.Range("I70").Value = 3993.5454
$ws.Range("K70").Value = 3993.5454
$ws.Range("M70").Value = -3723.5454
$ws.Range("H73").Value = 4284.9443
$ws.Range("I73").Value = 3993.5454
$ws.Range("K73").Value = 3993.5454
$ws.Range("M73").Value = -3057.5454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1506.1154
$ws.Range("I22").Value = 996
$ws.Range("J22").Value = 1627.5714
$ws.Range("K22").Value = 996
$ws.Range("L22").Value = 1627.5714
$ws.Range("M22").Value = -701
$ws.Range("N22").Value = -2217.5714
$ws.Range("H27").Value = 1506.1154
$ws.Range("I27").Value = 996
$ws.Range("J27").Value = 1627.5714
$ws.Range("K27").Value = 996
$ws.Range("L27").Value = 1627.5714
$ws.Range("M27").Value = -889
$ws.Range("N27").Value = -1841.5714
$ws.Range("H40").Value = 674866.5600000001
$ws.Range("I40").Value = 778422.9399999999
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 778422.9399999999
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -778286.9399999999
$ws.Range("N40").Value = -2022
$ws.Range("H46").Value = 1273.2593
$ws.Range("I46").Value = 1239.8334
$ws.Range("K46").Value = 1239.8334
$ws.Range("M46").Value = -1051.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2914.7
$ws.Range("I96").Value = 2215.1538
$ws.Range("J96").Value = 4213.857
$ws.Range("K96").Value = 2215.1538
$ws.Range("L96").Value = 4213.857
$ws.Range("M96").Value = -842.1538
$ws.Range("N96").Value = -6959.857
$ws.Range("H123").Value = 27854.455
$ws.Range("J123").Value = 27854.455
$ws.Range("L123").Value = 27854.455
$ws.Range("N123").Value = -37654.455
$ws.Range("H132").Value = 1709.9454
$ws.Range("I132").Value = 1182.6562
$ws.Range("J132").Value = 2443.5652
$ws.Range("K132").Value = 3547.9686
$ws.Range("L132").Value = 7330.6956
$ws.Range("M132").Value = -1017.9686
$ws.Range("N132").Value = -12390.6956
$ws.Range("H136").Value = 3633.2222
$ws.Range("I136").Value = 3826.5
$ws.Range("J136").Value = 3246.6667
$ws.Range("K136").Value = 11479.5
$ws.Range("L136").Value = 9740.000100000001
$ws.Range("M136").Value = -8929.5
$ws.Range("N136").Value = -14840.0001

Write-Output "Applied all changes"